$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.472.77'
$ws.Range("E2").Value = '  +0.13%  '

$ws.Range("D3").Value = '3.672.93'
$ws.Range("E3").Value = '  -0.44%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '645.60'
$ws.Range("E5").Value = '  -5.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.04'
$ws.Range("E6").Value = '  -0.26%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("E8").Value = '  +0.19%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.145'
$ws.Range("E9").Value = '  -0.99%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.04'
$ws.Range("E10").Value = '  -0.60%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.439'
$ws.Range("E11").Value = '  +0.31%  '

$ws.Range("E12").Value = '  -0.99%  '

$ws.Range("D13").Value = '4.294.58'
$ws.Range("E13").Value = '  -0.41%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.30'
$ws.Range("E14").Value = '  -0.55%  '

$ws.Range("D15").Value = '3.690.53'
$ws.Range("E15").Value = '  +0.25%  '

$ws.Range("D16").Value = '69.483.88'
$ws.Range("E16").Value = '  +0.21%  '

$ws.Range("E17").Value = '  +1.33%  '

$ws.Range("E18").Value = '  -0.40%  '

$ws.Range("E19").Value = '  +0.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '466.39'
$ws.Range("E20").Value = '  -0.50%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.72'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.643'
$ws.Range("E22").Value = '  -1.42%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.38'
$ws.Range("E23").Value = '  -0.67%  '

$ws.Range("D24").Value = '3.817.27'
$ws.Range("E24").Value = '  -0.47%  '

$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("E26").Value = '  +1.17%  '

$ws.Range("E27").Value = '  -1.36%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.89'
$ws.Range("E28").Value = '  -2.79%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.60'
$ws.Range("E29").Value = '  -2.81%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.65'
$ws.Range("E30").Value = '  -6.23%  '

$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("E32").Value = '  -0.11%  '

$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.58'
$ws.Range("E33").Value = '  -1.32%  '

$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.44'
$ws.Range("E34").Value = '  -2.79%  '

$ws.Range("D35").Value = '3.663.92'
$ws.Range("E35").Value = '  -0.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.161'
$ws.Range("E36").Value = '  +2.31%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.34'
$ws.Range("E37").Value = '  +0.78%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.90'
$ws.Range("E39").Value = '  -5.49%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '178.97'
$ws.Range("E40").Value = '  +4.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.18'
$ws.Range("E42").Value = '  -3.84%  '

$ws.Range("E43").Value = '  -1.49%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '46.71'
$ws.Range("E45").Value = '  -2.11%  '

$ws.Range("E46").Value = '  -0.63%  '

$ws.Range("E47").Value = '  -2.81%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.85'
$ws.Range("E48").Value = '  -5.17%  '

$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.79'
$ws.Range("E49").Value = '  -0.04%  '

$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000264'
$ws.Range("E50").Value = '  -4.83%  '

$ws.Range("E51").Value = '  -4.64%  '
